$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "63.672.31"
$ws.Range("E2").Value2 = "  +5.73%  "
$ws.Range("D3").Value2 = "3.431.80"
$ws.Range("E3").Value2 = "  +7.12%  "
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "577.77"
$ws.Range("E5").Value2 = "  +7.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "156.69"
$ws.Range("E6").Value2 = "  +7.15%  "
$ws.Range("E7").Value2 = "  +0.00%  "
$ws.Range("D8").Value2 = "3.437.18"
$ws.Range("E8").Value2 = "  +7.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.535"
$ws.Range("E9").Value2 = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.57"
$ws.Range("E10").Value2 = "  +3.05%  "
$ws.Range("E11").Value2 = "  +8.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.438"
$ws.Range("E12").Value2 = "  +0.99%  "
$ws.Range("D13").Value2 = "4.026.81"
$ws.Range("E13").Value2 = "  +7.31%  "
$ws.Range("E14").Value2 = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0000186"
$ws.Range("E15").Value2 = "  +7.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "27.25"
$ws.Range("E16").Value2 = "  +5.25%  "
$ws.Range("D17").Value2 = "63.812.82"
$ws.Range("E17").Value2 = "  +5.98%  "
$ws.Range("D18").Value2 = "3.433.71"
$ws.Range("E18").Value2 = "  +7.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.44"
$ws.Range("E19").Value2 = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "14.29"
$ws.Range("E20").Value2 = "  +7.45%  "
$ws.Range("E21").Value2 = "  +3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "392.25"
$ws.Range("E22").Value2 = "  +5.68%  "
$ws.Range("B23").Value2 = "Polygon"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.540"
$ws.Range("E23").Value2 = "  +3.09%  "
$ws.Range("B24").Value2 = "Dai"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.999"
$ws.Range("E24").Value2 = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "72.08"
$ws.Range("E25").Value2 = "  +3.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.0000108"
$ws.Range("E26").Value2 = "  +22.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.54"
$ws.Range("E27").Value2 = "  +10.11%  "
$ws.Range("E28").Value2 = "  +6.95%  "
$ws.Range("E29").Value2 = "  +0.04%  "
$ws.Range("E30").Value2 = "  +7.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.62"
$ws.Range("E31").Value2 = "  +7.71%  "
$ws.Range("E32").Value2 = "  +14.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.75"
$ws.Range("E33").Value2 = "  +8.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "23.47"
$ws.Range("E34").Value2 = "  +4.42%  "
$ws.Range("E35").Value2 = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "6.80"
$ws.Range("E36").Value2 = "  +3.43%  "
$ws.Range("E37").Value2 = "  +9.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "158.88"
$ws.Range("E38").Value2 = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "28.06"
$ws.Range("E39").Value2 = "  +5.72%  "
$ws.Range("E40").Value2 = "  +10.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.89"
$ws.Range("E41").Value2 = "  +11.16%  "
$ws.Range("D42").Value2 = "2.925.12"
$ws.Range("E42").Value2 = "  +4.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0321"
$ws.Range("E43").Value2 = "  +2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.769"
$ws.Range("E44").Value2 = "  +6.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "41.82"
$ws.Range("E45").Value2 = "  +4.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "4.35"
$ws.Range("E46").Value2 = "  +3.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.09"
$ws.Range("E47").Value2 = "  +10.20%  "
$ws.Range("D48").Value2 = "3.483.84"
$ws.Range("E48").Value2 = "  +7.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "22.64"
$ws.Range("E50").Value2 = "  +3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "296.39"
$ws.Range("E51").Value2 = "  +12.39%  "
